# Daily "water delivery countdown" roll-forward:
#   - Column D = total days for the current delivery cycle
#   - Column E = days remaining in the current cycle
#   - Column F = cycle start date (stored as a plain yyyyMMdd number)
#
# One day has elapsed since the last update, so for every data row:
#   - if the remaining-days counter was about to hit zero (E = 1), the
#     delivery is renewed: E resets to the full cycle length (D) and the
#     start date F jumps forward by 10 (the next restock offset)
#   - otherwise E simply ticks down by 1 and F is left untouched
#
# Rows whose start date isn't a well-formed 8-digit yyyyMMdd value (a bad
# data entry) are left completely untouched, same as upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$firstRow = $usedRange.Row
$lastRow = $firstRow + $usedRange.Rows.Count - 1

for ($r = $firstRow + 1; $r -le $lastRow; $r++) {

    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $fCell = $ws.Cells.Item($r, 6)

    $dVal = $dCell.Value()
    $eVal = $eCell.Value()
    $fVal = $fCell.Value()

    if ($eVal -eq $null -or $dVal -eq $null -or $fVal -eq $null) {
        continue
    }

    $fStr = [string]$fVal
    if ($fStr.Length -ne 8) {
        # malformed start-date (e.g. row 36's "202510929") - skip untouched
        continue
    }

    if ($eVal -eq 1) {
        $eCell.Value = $dVal
        $fCell.Value = $fVal + 10
    } else {
        $eCell.Value = $eVal - 1
    }
}
